$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44304, 0, 3, 112.4016485575122),
    @(44305, 0, 1, 37.46721618583739),
    @(44306, 0, 0, 0),
    @(44307, 0, 0, 0)
)

$startRow = 230
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the format of column A from the row above so the new date cell
    # picks up the same style (s="2": centered, bordered date format).
    $ws.Range("A$($r - 1)").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
